function Set-TextCell($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the 2022-Q1 summary row to the "总计" sheet FIRST (before inserting
#    any new sheets, since adding a sheet shifts the position of "总计" and
#    any worksheet reference grabbed beforehand would otherwise go stale).
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Cells.Item(3, 1).Copy($totalSheet.Cells.Item(2, 1))

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 19
$totalSheet.Cells.Item(2, 4).Value = 3.77

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(7, 1).Value = 5

# ---------------------------------------------------------------------------
# 2. Insert the new "2022-Q1" sheet right before the "总计" summary sheet,
#    using an existing quarter sheet as the style template so fonts/borders
#    match the rest of the workbook. Re-fetch "总计" by name right before
#    use, since the earlier reference is now positionally stale once a new
#    sheet gets inserted ahead of it.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newSheet.Name = "2022-Q1"

# Header row (identical wording/style to the other quarter sheets).
$template.Range("B1:H1").Copy($newSheet.Range("B1:H1"))

# Style template for the "index" column (A) used by every data row.
$indexTemplate = $template.Cells.Item(2, 1)

# Fund holdings data, columns kept as parallel flat arrays (this PS engine
# does not support nested @() array literals).
$codes = @("519732", "007449", "398001", "010861", "006348", "006642", "970023", "005589", "005904", "007450", "007592", "002732", "000354", "005905", "008116", "012297", "009327", "012298", "009328")
$names = @("交银定期支付双息平衡混合", "兴全多维价值混合A", "中海优质成长混合", "长信企业优选一年持有期灵活配置混合", "银华盛利混合", "华泰保兴吉年利定期开放混合", "天风天盈一年定期开放混合", "长信企业精选两年定期开放灵活配置混合", "华泰保兴成长优选混合A", "兴全多维价值混合C", "华夏价值精选混合", "长盛沪港深优势精选灵活配置混合", "长盛城镇化主题混合", "华泰保兴成长优选混合C", "银华沪深股通精选混合", "东兴宸瑞量化混合型证券投资基金A", "东兴兴晟混合A", "东兴宸瑞量化混合型证券投资基金C", "东兴兴晟混合C")
$scales = @("40.83", "28.03", "14.42", "9.39", "8.41", "7.05", "2.97", "5.84", "4.77", "4.48", "2.55", "0.74", "0.43", "0.46", "0.55", "0.48", "0.34", "0.21", "0.08")
$stockPos = @("67.67", "84.40", "90.86", "80.21", "88.43", "91.97", "39.47", "79.99", "82.40", "84.40", "94.58", "76.85", "78.41", "82.40", "88.86", "89.42", "79.83", "89.42", "79.83")
$posShare = @("2.03", "2.60", "3.85", "2.93", "3.26", "3.54", "7.03", "2.92", "3.14", "2.60", "4.56", "3.87", "4.71", "3.14", "2.60", "1.89", "1.26", "1.89", "1.26")
$heldValue = @("0.8288", "0.7288", "0.5552", "0.2751", "0.2742", "0.2496", "0.2088", "0.1705", "0.1498", "0.1165", "0.1163", "0.0286", "0.0203", "0.0144", "0.0143", "0.0091", "0.0043", "0.0040", "0.0010")
$posRank = @(8, 6, 7, 10, 5, 9, 1, 7, 9, 6, 6, 6, 8, 9, 8, 2, 2, 2, 2)

for ($i = 0; $i -lt $codes.Count; $i++) {
    $r = $i + 2

    $indexTemplate.Copy($newSheet.Cells.Item($r, 1))
    $newSheet.Cells.Item($r, 1).Value = $i

    Set-TextCell $newSheet.Cells.Item($r, 2) $codes[$i]
    Set-TextCell $newSheet.Cells.Item($r, 3) $names[$i]
    Set-TextCell $newSheet.Cells.Item($r, 4) $scales[$i]
    Set-TextCell $newSheet.Cells.Item($r, 5) $stockPos[$i]
    Set-TextCell $newSheet.Cells.Item($r, 6) $posShare[$i]
    Set-TextCell $newSheet.Cells.Item($r, 7) $heldValue[$i]
    $newSheet.Cells.Item($r, 8).Value = $posRank[$i]
}

Write-Output "done"
